$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.410.74'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.478.75'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.47%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.81'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.72%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.549'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.10%  '
$ws.Range('E9').Value = '  -4.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.62'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.24%  '
$ws.Range('E11').Value = '  -2.66%  '
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.01'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.859.68'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.496.78'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.22'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('E17').Value = '  -3.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.393.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.33'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.09%  '
$ws.Range('E20').Value = '  -2.58%  '
$ws.Range('E21').Value = '  -8.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.81'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.76'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.19%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.14'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.72'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.62'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '152.50'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.43%  '
$ws.Range('E32').Value = '  -6.54%  '
$ws.Range('E33').Value = '  -3.07%  '
$ws.Range('E34').Value = '  -7.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0749'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.08'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.65'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.02%  '
$ws.Range('E38').Value = '  -4.69%  '
$ws.Range('E39').Value = '  -2.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.27'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.11%  '
$ws.Range('E41').Value = '  -8.31%  '
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('E43').Value = '  -10.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.986.39'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.70%  '
$ws.Range('E45').Value = '  -4.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.75'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.723.45'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '69.77'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '97.20'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '74.65'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.64%  '
